# Update column F ("dSF") values for several rows on the active sheet,
# per the re-pulled data / mean calculation described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -7
    4  = -5
    5  = -1
    6  = -6
    7  = -2
    9  = -3
    10 = 13
    13 = -3
    16 = 0
    18 = -2
    20 = -6
    22 = -2
    23 = -4
    24 = -6
    25 = 8
    26 = 2
    28 = -5
    30 = -3
    31 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
